# Remove the "Dancer in the Darkness" stat-block row (row 14).
# This also collapses the table by one row, moving the "Krashtkid"
# row up from row 15 to row 14, and drops the now-unused shared
# strings ("Dancer in the Darkness", "3d6x2", "Earth") from the
# workbook when it is saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(14).Delete()

# Mirror Excel's natural post-delete selection behaviour: the
# selection lands on the row that shifted up into the deleted row's
# place (now the last row of data).
$ws.Range("A14:XFD14").Select()
